$d = $word.ActiveDocument

$replacements = @(
    @{old = "A : 15"; new = "A : 12"},
    @{old = "A-: 14"; new = "A-: 11"},
    @{old = "B+: 13"; new = "B+: 10"},
    @{old = "B : 12"; new = "B : 9"},
    @{old = "B-: 11"; new = "B-: 8"},
    @{old = "C+: 10"; new = "C+: 7"},
    @{old = "C : 9"; new = "C : 6"},
    @{old = "C-: 7"; new = "C-: 5"},
    @{old = "P : 7"; new = "P : 5"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
